# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output (commit: "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 4185
$ws1.Range("F4").Value  = 2394
$ws1.Range("F7").Value  = 41
$ws1.Range("F8").Value  = 43
$ws1.Range("F9").Value  = 211
$ws1.Range("F10").Value = 122
$ws1.Range("F11").Value = 115
$ws1.Range("F12").Value = 148
$ws1.Range("F13").Value = 1564
$ws1.Range("F14").Value = 290
$ws1.Range("F15").Value = 3165
$ws1.Range("F16").Value = 213

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 4185
$ws4.Range("F4").Value  = 2394
$ws4.Range("F8").Value  = 41
$ws4.Range("F9").Value  = 43
$ws4.Range("F11").Value = 211
$ws4.Range("F12").Value = 122
$ws4.Range("F13").Value = 115
$ws4.Range("F14").Value = 148
$ws4.Range("F17").Value = 1564
$ws4.Range("F18").Value = 290
$ws4.Range("F19").Value = 3166
$ws4.Range("F20").Value = 213
